# Adds a new "2022-Q1" sheet (repurposing the existing "总计" sheet, which
# is recreated fresh at the end) and updates the "总计" totals sheet with a
# new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: force a cell to be stored as literal text (even if its content
# looks numeric, e.g. fund codes with leading zeros, or decimal strings
# like "6.33") without leaving a stray NumberFormat style behind.
# ---------------------------------------------------------------------
function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# =======================================================================
# Step 1: the current "总计" sheet (5th sheet) becomes "2022-Q1" and is
# repopulated with the quarterly fund-holding detail data.
# =======================================================================
$q1 = $wb.Worksheets.Item(5)

# Copy the existing header/index-column formatting (style index 2) onto
# the newly-used cells (E1:H1 header cells, and A6:A11 row-index cells).
$fmtSrcHeader = $q1.Range("B1")
$fmtSrcHeader.Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

$fmtSrcIndex = $q1.Range("A2")
$fmtSrcIndex.Copy()
$q1.Range("A6:A11").PasteSpecial(-4122)

$q1.Application.CutCopyMode = $false

# Header row
$q1.Cells.Item(1,2).Value = "基金代码"
$q1.Cells.Item(1,3).Value = "基金名称"
$q1.Cells.Item(1,4).Value = "基金规模"
$q1.Cells.Item(1,5).Value = "股票总仓位"
$q1.Cells.Item(1,6).Value = "仓位占比"
$q1.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1.Cells.Item(1,8).Value = "仓位排名"

# Data rows: fund code, fund name, fund size, stock position, position
# ratio, held market value, position rank.
$q1Rows = @(
    @("001325", "鹏华弘和灵活配置混合A", "6.33", "22.78", "1.72", "0.1089", 3),
    @("001331", "鹏华弘信灵活配置混合A", "5.01", "25.87", "1.87", "0.0937", 3),
    @("003142", "鹏华弘达灵活配置混合A", "4.57", "24.46", "1.54", "0.0704", 8),
    @("003780", "鹏华兴悦定期开放灵活配置混合", "4.34", "20.74", "1.62", "0.0703", 5),
    @("003663", "鹏华兴泰定期开放灵活配置混合", "4.70", "22.85", "1.49", "0.0700", 7),
    @("001327", "鹏华弘华灵活配置混合A", "2.95", "28.77", "1.85", "0.0546", 4),
    @("001326", "鹏华弘和灵活配置混合C", "2.40", "22.78", "1.72", "0.0413", 3),
    @("001328", "鹏华弘华灵活配置混合C", "0.99", "28.77", "1.85", "0.0183", 4),
    @("001332", "鹏华弘信灵活配置混合C", "0.87", "25.87", "1.87", "0.0163", 3),
    @("003143", "鹏华弘达灵活配置混合C", "0.62", "24.46", "1.54", "0.0095", 8)
)

for ($i = 0; $i -lt $q1Rows.Count; $i++) {
    $r = $i + 2
    $rowData = $q1Rows[$i]

    $q1.Cells.Item($r, 1).Value = $i

    Set-TextValue $q1.Cells.Item($r, 2) $rowData[0]
    Set-TextValue $q1.Cells.Item($r, 3) $rowData[1]
    Set-TextValue $q1.Cells.Item($r, 4) $rowData[2]
    Set-TextValue $q1.Cells.Item($r, 5) $rowData[3]
    Set-TextValue $q1.Cells.Item($r, 6) $rowData[4]
    Set-TextValue $q1.Cells.Item($r, 7) $rowData[5]

    $q1.Cells.Item($r, 8).Value = $rowData[6]
}

$q1.Name = "2022-Q1"

# =======================================================================
# Step 2: append a brand-new sheet at the end of the workbook and name it
# "总计"; populate it with the historical totals plus the new 2022-Q1 row.
# =======================================================================
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$totals = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$totals.Name = "总计"

# Bring over the same header/index-column styling (style index 2) used on
# the other sheets, sourced from the "2022-Q1" sheet we just built.
$fmtSrcHeader.Copy()
$totals.Range("B1:D1").PasteSpecial(-4122)

$fmtSrcIndex.Copy()
$totals.Range("A2:A6").PasteSpecial(-4122)

$totals.Application.CutCopyMode = $false

$totals.Cells.Item(1,2).Value = "日期"
$totals.Cells.Item(1,3).Value = "持有数量(只)"
$totals.Cells.Item(1,4).Value = "持有市值(亿元)"

$totalsRows = @(
    @("2022-Q1", 10, 0.55),
    @("2021-Q4", 1, 0.07000000000000001),
    @("2021-Q3", 1, 0.09),
    @("2021-Q2", 3, 0.15),
    @("2021-Q1", 10, 0.85)
)

for ($i = 0; $i -lt $totalsRows.Count; $i++) {
    $r = $i + 2
    $rowData = $totalsRows[$i]

    $totals.Cells.Item($r, 1).Value = $i
    $totals.Cells.Item($r, 2).Value = $rowData[0]
    $totals.Cells.Item($r, 3).Value = $rowData[1]
    $totals.Cells.Item($r, 4).Value = $rowData[2]
}

# Restore the originally-active sheet/tab selection (the first sheet),
# since adding/activating new sheets above shifted it.
$wb.Worksheets.Item(1).Activate()
